$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1226.6666
$ws.Range("I34").Value = 1226.6666
$ws.Range("K34").Value = 1226.6666
$ws.Range("M34").Value = -1023.6666
$ws.Range("H36").Value = 1226.6666
$ws.Range("I36").Value = 1226.6666
$ws.Range("K36").Value = 1226.6666
$ws.Range("M36").Value = -511.6666
$ws.Range("H43").Value = 14200.375
$ws.Range("J43").Value = 14800.429
$ws.Range("L43").Value = 14800.429
$ws.Range("N43").Value = -14938.429
$ws.Range("H64").Value = 3319.75
$ws.Range("I64").Value = 3341.6667
$ws.Range("K64").Value = 3341.6667
$ws.Range("M64").Value = -3093.6667
$ws.Range("H67").Value = 3319.75
$ws.Range("I67").Value = 3341.6667
$ws.Range("K67").Value = 3341.6667
$ws.Range("M67").Value = -2483.6667
$ws.Range("H74").Value = 2971.7693
$ws.Range("I74").Value = 2535.96
$ws.Range("K74").Value = 2535.96
$ws.Range("M74").Value = -1599.96
$ws.Range("H77").Value = 2971.7693
$ws.Range("I77").Value = 2535.96
$ws.Range("K77").Value = 12679.8
$ws.Range("M77").Value = -7999.799999999999
$ws.Range("H87").Value = 19354
$ws.Range("J87").Value = 19354
$ws.Range("L87").Value = 19354
$ws.Range("N87").Value = -21850
$ws.Range("H90").Value = 19354
$ws.Range("J90").Value = 19354
$ws.Range("L90").Value = 58062
$ws.Range("N90").Value = -70542

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11629434
$ws.Range("I61").Value = 13514971
$ws.Range("K61").Value = 13514971
$ws.Range("M61").Value = -13514759
$ws.Range("H110").Value = 1373.75
$ws.Range("I110").Value = 1220
$ws.Range("J110").Value = 1630
$ws.Range("K110").Value = 1220
$ws.Range("L110").Value = 1630
$ws.Range("M110").Value = 825
$ws.Range("N110").Value = -5720
$ws.Range("H136").Value = 11629434
$ws.Range("I136").Value = 13514971
$ws.Range("K136").Value = 40544913
$ws.Range("M136").Value = -40542363

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2131.0588
$ws.Range("I20").Value = 2015.6428
$ws.Range("J20").Value = 2669.6667
$ws.Range("K20").Value = 2015.6428
$ws.Range("L20").Value = 2669.6667
$ws.Range("M20").Value = -1768.6428
$ws.Range("N20").Value = -3163.6667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 11156.3
$ws.Range("I26").Value = 1642.8572
$ws.Range("K26").Value = 1642.8572
$ws.Range("M26").Value = -1355.8572
$ws.Range("H44").Value = 123333.336
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30884
$ws.Range("H54").Value = 29055.2
$ws.Range("J54").Value = 29055.2
$ws.Range("L54").Value = 29055.2
$ws.Range("N54").Value = -30371.2
$ws.Range("H56").Value = 31449.5
$ws.Range("J56").Value = 31449.5
$ws.Range("L56").Value = 31449.5
$ws.Range("N56").Value = -33139.5
$ws.Range("H140").Value = 43006.555
$ws.Range("J140").Value = 43006.555
$ws.Range("L140").Value = 43006.555
$ws.Range("N140").Value = -53366.555
$ws.Range("H12").Value = 12500281
$ws.Range("I12").Value = 12500281
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 12500281
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -12500111
$ws.Range("N12").ClearContents()
$ws.Range("H105").Value = 1636.125
$ws.Range("I105").Value = 1636.125
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1636.125
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 110.875
$ws.Range("N105").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6876.136
$ws.Range("I3").Value = 3061.3635
$ws.Range("J3").Value = 10690.909
$ws.Range("K3").Value = 9184.0905
$ws.Range("L3").Value = 32072.727
$ws.Range("M3").Value = -9072.0905
$ws.Range("N3").Value = -32296.727
$ws.Range("H4").Value = 176.7619
$ws.Range("I4").Value = 176
$ws.Range("J4").Value = 192
$ws.Range("K4").Value = 528
$ws.Range("L4").Value = 576
$ws.Range("M4").Value = -416
$ws.Range("N4").Value = -800
$ws.Range("H7").Value = 14285858
$ws.Range("I7").Value = 25000118
$ws.Range("J7").Value = 179
$ws.Range("K7").Value = 75000354
$ws.Range("L7").Value = 537
$ws.Range("M7").Value = -75000242
$ws.Range("N7").Value = -761
$ws.Range("H25").Value = 886.8889
$ws.Range("I25").Value = 260
$ws.Range("J25").Value = 1200.3334
$ws.Range("K25").Value = 780
$ws.Range("L25").Value = 3601.0002
$ws.Range("M25").Value = -611
$ws.Range("N25").Value = -3939.0002
$ws.Range("H30").Value = 886.8889
$ws.Range("I30").Value = 260
$ws.Range("J30").Value = 1200.3334
$ws.Range("K30").Value = 780
$ws.Range("L30").Value = 3601.0002
$ws.Range("M30").Value = -678
$ws.Range("N30").Value = -3805.0002
$ws.Range("H34").Value = 1299.3334
$ws.Range("I34").Value = 247.5
$ws.Range("J34").Value = 1599.8572
$ws.Range("K34").Value = 742.5
$ws.Range("L34").Value = 4799.571599999999
$ws.Range("M34").Value = -658.5
$ws.Range("N34").Value = -4967.571599999999
$ws.Range("H37").Value = 91110
$ws.Range("J37").Value = 91110
$ws.Range("L37").Value = 273330
$ws.Range("N37").Value = -273554
$ws.Range("H39").Value = 494.64285
$ws.Range("J39").Value = 499.39026
$ws.Range("L39").Value = 1498.17078
$ws.Range("N39").Value = -2086.17078
$ws.Range("H109").Value = 3615.625
$ws.Range("I109").Value = 1370.25
$ws.Range("J109").Value = 4064.7
$ws.Range("K109").Value = 4110.75
$ws.Range("L109").Value = 12194.1
$ws.Range("M109").Value = -3070.75
$ws.Range("N109").Value = -14274.1
$ws.Range("H118").Value = 1073.2222
$ws.Range("I118").Value = 200
$ws.Range("K118").Value = 600
$ws.Range("M118").Value = 643
$ws.Range("H131").Value = 843.86
$ws.Range("J131").Value = 859.44794
$ws.Range("L131").Value = 2578.34382
$ws.Range("N131").Value = -12658.34382
$ws.Range("H134").Value = 3864.4443
$ws.Range("I134").Value = 2190.5881
$ws.Range("J134").Value = 6710
$ws.Range("K134").Value = 6571.7643
$ws.Range("L134").Value = 20130
$ws.Range("M134").Value = -1501.7643
$ws.Range("N134").Value = -30270
$ws.Range("H137").Value = 6671173
$ws.Range("I137").Value = 13890584
$ws.Range("J137").Value = 7102.231
$ws.Range("K137").Value = 41671752
$ws.Range("L137").Value = 21306.693
$ws.Range("M137").Value = -41666652
$ws.Range("N137").Value = -31506.693
$ws.Range("H103").Value = 7944.1177
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 7944.1177
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 23832.3531
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -25590.3531

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5138.8887
$ws.Range("J70").Value = 5138.8887
$ws.Range("L70").Value = 5138.8887
$ws.Range("N70").Value = -5678.8887
$ws.Range("H73").Value = 5138.8887
$ws.Range("J73").Value = 5138.8887
$ws.Range("L73").Value = 5138.8887
$ws.Range("N73").Value = -7010.8887

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 37224.5
$ws.Range("J42").Value = 37224.5
$ws.Range("L42").Value = 37224.5
$ws.Range("N42").Value = -37980.5
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108
